$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-07-20 Sunday" "2025-07-21 Monday"

Replace-Text "173÷7=" "408÷8="
Replace-Text "328÷2=" "324÷2="
Replace-Text "941÷9=" "634÷7="
Replace-Text "880÷6=" "716÷7="
Replace-Text "643÷7=" "513÷4="

Replace-Text "726÷7=" "181÷4="
Replace-Text "227÷5=" "615÷9="
Replace-Text "107÷3=" "708÷4="
Replace-Text "986÷8=" "480÷7="
Replace-Text "930÷5=" "512÷7="

Replace-Text "824÷7=" "883÷7="
Replace-Text "827÷6=" "503÷4="
Replace-Text "753÷4=" "281÷6="
Replace-Text "497÷2=" "302÷5="
Replace-Text "220÷8=" "473÷6="

Replace-Text "269÷4=" "345÷7="
Replace-Text "681÷5=" "916÷7="
Replace-Text "364÷2=" "376÷7="
Replace-Text "804÷3=" "242÷2="
Replace-Text "453÷4=" "947÷5="

Replace-Text "175÷5=" "732÷7="
Replace-Text "437÷9=" "300÷7="
Replace-Text "448÷3=" "916÷7="
Replace-Text "425÷3=" "524÷7="
Replace-Text "509÷9=" "848÷7="
